# Added a tester set timeout.
# Both the "GET_Tests" and "POST Tests" sheets get a new column L that
# mirrors the existing J/K "timeout" config columns (value 2) for every
# data row, the selection is moved onto the freshly-added column, and the
# POST Tests sheet's wrapped-text rows grow taller to match their new
# (wider) layout.

$wb = $excel.ActiveWorkbook

$wsGet  = $wb.Worksheets.Item("GET_Tests")
$wsPost = $wb.Worksheets.Item("POST Tests")

# --- GET_Tests: add the new timeout column (L) ---------------------------
$wsGet.Range("L1").Value = 2
$wsGet.Range("L2").Value = 2
$wsGet.Range("L3").Value = 2

# --- POST Tests: add the new timeout column (L) ---------------------------
$wsPost.Range("L1").Value = 2
$wsPost.Range("L2").Value = 2
$wsPost.Range("L3").Value = 2

# The wrapped "headers" column on POST Tests reflows with the sheet's new
# width, so its data rows grow taller.
$wsPost.Rows.Item(2).RowHeight = 186.35
$wsPost.Rows.Item(3).RowHeight = 186.35

# --- Update each sheet's remembered selection to the new column ----------
$wsPost.Activate()
$wsPost.Range("L1:L3").Select()

$wsGet.Activate()
$wsGet.Range("L1:L3").Select()
